# Update res_bus vm_pu results for "case with 380 kV" (Case_4_163).
# Columns updated per row: B,C,D,E,F,I,J,K,L,M,N (G stays 1, H stays blank, A is the index column).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @(
    @(1.02, 1.056882949518213, 1.062984988745737, 1.060571346374045, 1.070800210076821, 1.05257268618011, 1.061882265411271, 1.065704854516345, 1.063297763419177, 1.073499082374265, 1.024126665081736),
    @(1.02, 1.058146015546485, 1.064038772042396, 1.061785088383639, 1.072042952165762, 1.053033773554781, 1.062796071898278, 1.066572960422947, 1.064324949209973, 1.074557201794975, 1.024456097618107),
    @(1.02, 1.058962553132144, 1.064719982682444, 1.062570045358446, 1.072846797705581, 1.053330450835493, 1.063386098300243, 1.067133423236317, 1.064988623400094, 1.075241018051614, 1.024668250918759),
    @(1.02, 1.059305649172304, 1.065006208325222, 1.062899944524994, 1.073184667163079, 1.053454773765741, 1.063633844222782, 1.067368742265766, 1.065267398958143, 1.075528291438865, 1.024757198699166),
    @(1.02, 1.05936324626201, 1.065054257824256, 1.06295533045096, 1.073241393066767, 1.053475624702616, 1.063675424257891, 1.067408235860758, 1.065314193004077, 1.075576514071718, 1.024772119283664),
    @(1.02, 1.058967138287039, 1.064723807851105, 1.062574453864064, 1.07285131259461, 1.053332113615463, 1.06338940987398, 1.067136568755126, 1.064992349325526, 1.075244857407224, 1.024669440391265),
    @(1.019999999999999, 1.057309966019477, 1.063341257189803, 1.0609816237616, 1.071220262562593, 1.052728860850925, 1.062191354401369, 1.065998497049375, 1.063645110812452, 1.073856857828595, 1.024238208229824),
    @(1.02, 1.054383918583716, 1.060899894072918, 1.058171553208824, 1.068343792060776, 1.051652949551294, 1.060070418267762, 1.063983328817272, 1.061263461292778, 1.071404348657589, 1.023470542761449),
    @(1.02, 1.052429043714747, 1.059268731988833, 1.056295783883554, 1.066424396504301, 1.05092691725153, 1.058649721402857, 1.062633203442044, 1.059670401508755, 1.069764701652298, 1.022953489979003),
    @(1.02, 1.051581523134824, 1.058561539081429, 1.05548294991707, 1.065592824366862, 1.050610441069906, 1.058032912781672, 1.06204697011391, 1.058979298455789, 1.06905358115719, 1.022728338644489),
    @(1.02, 1.051266554790579, 1.058298719825393, 1.055180931613549, 1.065283868797759, 1.050492570617447, 1.057803553674781, 1.061828970894364, 1.058722393942426, 1.06878926462798, 1.02264451655527),
    @(1.02, 1.051334123943628, 1.058355101661421, 1.055245719944856, 1.065350144157893, 1.050517868606705, 1.057852763292794, 1.061875743644411, 1.058777509845712, 1.06884596937629, 1.02266250532995),
    @(1.02, 1.051555491071338, 1.058539817156723, 1.055457986959401, 1.065567287504708, 1.050600704340264, 1.058013958981118, 1.062028955265341, 1.058958066707138, 1.069031736252325, 1.022721413783393),
    @(1.02, 1.051691861143592, 1.058653608322733, 1.055588758905635, 1.065701066959461, 1.050651700098976, 1.058113243893192, 1.062123321316793, 1.059069287409381, 1.069146170112467, 1.022757683906774),
    @(1.02, 1.052485268451306, 1.059315647087842, 1.056349715781071, 1.06647957520486, 1.05094787635679, 1.05869062221162, 1.062672075456838, 1.059716240137031, 1.069811872044725, 1.022968405818361),
    @(1.02, 1.052982668668664, 1.059730686565277, 1.056826877278096, 1.066967787014748, 1.051133096632916, 1.059052355922763, 1.063015858806167, 1.060121706879841, 1.070229141234517, 1.023100247011876),
    @(1.02, 1.053272693077699, 1.059972686487281, 1.057105138547509, 1.067252508433853, 1.051240930039614, 1.05926319127879, 1.063216225678981, 1.060358083777881, 1.070472417265112, 1.023177025849811),
    @(1.02, 1.053371566947992, 1.060055187773504, 1.057200008551525, 1.067349583672041, 1.051277664157136, 1.059335054008861, 1.063284519200602, 1.060438661084104, 1.070555349534125, 1.023203184801536),
    @(1.02, 1.052929312782541, 1.059686165630631, 1.056775688489587, 1.066915411112544, 1.051113245203114, 1.059013561662791, 1.062978990276645, 1.060078217087892, 1.070184383603486, 1.023086114315462),
    @(1.02, 1.05149030846058, 1.058485426832721, 1.055395482280084, 1.06550334623384, 1.050576320048728, 1.057966497763449, 1.061983845075487, 1.058904902704404, 1.068977037385103, 1.022704071988364),
    @(1.02, 1.05058461305844, 1.057729684593898, 1.054527137178132, 1.064615102985034, 1.050236898362967, 1.057306725501409, 1.061356732844065, 1.05816604587213, 1.068216918526101, 1.022462762093199),
    @(1.02, 1.051064829476028, 1.058130393672549, 1.054987517064914, 1.065086018665127, 1.050417006775263, 1.057656620883193, 1.061689312859987, 1.058557837670464, 1.068619968902054, 1.022590790065852),
    @(1.02, 1.052953422321852, 1.059706282989005, 1.056798818681802, 1.066939077664963, 1.051122215835219, 1.05903109160647, 1.062995650060517, 1.060097868638993, 1.070204607990267, 1.023092500649283),
    @(1.02, 1.055141092623535, 1.061531667117504, 1.058898432373278, 1.069087724526601, 1.05193263598572, 1.060619908363563, 1.064505465850571, 1.061880096413079, 1.07203918787609, 1.023669929241966)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range($cols[$j] + $rowNum).Value = $rowVals[$j]
    }
}
